$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.455.12"
$ws.Range("E2").Value = "  +9.09%  "
$ws.Range("D3").Value = "1.603.31"
$ws.Range("E3").Value = "  +8.35%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.71%  "
$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").Value = "0.9912"
$ws.Range("E5").Value = "  +2.07%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "303.54"
$ws.Range("E6").Value = "  +8.70%  "
$ws.Range("D7").Value = "0.3688"
$ws.Range("E7").Value = "  +0.73%  "
$ws.Range("D8").Value = "0.3405"
$ws.Range("E8").Value = "  +10.59%  "
$ws.Range("D9").Value = "42.69"
$ws.Range("E9").Value = "  +6.72%  "
$ws.Range("D10").Value = "1.143"
$ws.Range("E10").Value = "  +7.42%  "
$ws.Range("D11").Value = "0.07052"
$ws.Range("E11").Value = "  +5.62%  "
$ws.Range("D12").Value = "0.9991"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").Value = "5.939"
$ws.Range("E13").Value = "  +7.46%  "
$ws.Range("D14").Value = "19.71"
$ws.Range("E14").Value = "  +9.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.630"
$ws.Range("E15").Value = "  +6.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001085"
$ws.Range("E16").Value = "  +5.47%  "
$ws.Range("D17").Value = "1.599.15"
$ws.Range("E17").Value = "  +8.12%  "
$ws.Range("D18").Value = "0.9899"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06790"
$ws.Range("E19").Value = "  +14.32%  "
$ws.Range("D20").Value = "77.88"
$ws.Range("E20").Value = "  +11.75%  "
$ws.Range("D21").Value = "16.12"
$ws.Range("E21").Value = "  +11.04%  "
$ws.Range("D22").Value = "6.019"
$ws.Range("E22").Value = "  +9.44%  "
$ws.Range("D23").Value = "11.82"
$ws.Range("E23").Value = "  +6.94%  "
$ws.Range("D24").Value = "22.451.27"
$ws.Range("E24").Value = "  +8.76%  "
$ws.Range("E25").Value = "  +6.14%  "
$ws.Range("D26").Value = "2.522"
$ws.Range("E26").Value = "  +18.39%  "
$ws.Range("D27").Value = "150.69"
$ws.Range("E27").Value = "  +6.16%  "
$ws.Range("D28").Value = "19.52"
$ws.Range("D29").Value = "1.780.24"
$ws.Range("E29").Value = "  +8.61%  "
$ws.Range("D30").Value = "120.83"
$ws.Range("E30").Value = "  +5.86%  "
$ws.Range("D31").Value = "4.197"
$ws.Range("E31").Value = "  +6.98%  "
$ws.Range("D32").Value = "6.061"
$ws.Range("E32").Value = "  +20.78%  "
$ws.Range("D33").Value = "0.9525"
$ws.Range("E33").Value = "  +15.59%  "
$ws.Range("D34").Value = "0.08276"
$ws.Range("E34").Value = "  +3.60%  "
$ws.Range("D35").Value = "1.636"
$ws.Range("E35").Value = "  +6.65%  "
$ws.Range("D36").Value = "5.295"
$ws.Range("E36").Value = "  +11.88%  "
$ws.Range("D37").Value = "1.271"
$ws.Range("E37").Value = "  +5.24%  "
$ws.Range("D38").Value = "11.87"
$ws.Range("E38").Value = "  +13.23%  "
$ws.Range("D39").Value = "8.624"
$ws.Range("E39").Value = "  +12.62%  "
$ws.Range("D40").Value = "0.06116"
$ws.Range("E40").Value = "  +5.63%  "
$ws.Range("D41").Value = "0.02221"
$ws.Range("E41").Value = "  +8.51%  "
$ws.Range("D42").Value = "0.2028"
$ws.Range("E42").Value = "  +7.95%  "
$ws.Range("D43").Value = "0.9905"
$ws.Range("E43").Value = "  +1.95%  "
$ws.Range("D44").Value = "0.5918"
$ws.Range("E44").Value = "  +11.57%  "
$ws.Range("D45").Value = "3.851"
$ws.Range("E45").Value = "  +8.99%  "
$ws.Range("D46").Value = "13.29"
$ws.Range("E46").Value = "  +8.29%  "
$ws.Range("D47").Value = "0.5696"
$ws.Range("E47").Value = "  +9.58%  "
$ws.Range("D48").Value = "126.94"
$ws.Range("E48").Value = "  +7.13%  "
$ws.Range("D49").Value = "1.969"
$ws.Range("E49").Value = "  +8.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06820"
$ws.Range("E50").Value = "  +5.14%  "
$ws.Range("D51").Value = "73.82"
$ws.Range("E51").Value = "  +9.15%  "
